$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the post entry "「希望の窓は開いている」" which occupied row 89.
# Deleting the entire row shifts all subsequent rows (90-186) up by one
# (89-185) and Excel automatically updates the sheet dimension.
$ws.Rows.Item(89).Delete()
